# "how did I miss TAB()"
# Mark additional cells as supported ("X") and add a couple of footnotes
# to the Supported.xlsx compatibility matrix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CHDIR - Token/Parse/Eval now supported
$ws.Range("C12").Value = "X"
$ws.Range("D12").Value = "X"
$ws.Range("E12").Value = "X"

# COMMON - Eval now supported
$ws.Range("E21").Value = "X"

# ERROR - Token/Parse/Eval now supported
$ws.Range("C50").Value = "X"
$ws.Range("D50").Value = "X"
$ws.Range("E50").Value = "X"

# KEY - Token/Parse/Eval now supported, but not every parameter combo
$ws.Range("C72").Value = "X"
$ws.Range("D72").Value = "X"
$ws.Range("E72").Value = "X"
$ws.Range("G72").Value = "Not all params"

# RESUME - Token/Parse/Eval now supported
$ws.Range("C141").Value = "X"
$ws.Range("D141").Value = "X"
$ws.Range("E141").Value = "X"

# SCREEN (statement) - note about border color argument
$ws.Range("G150").Value = "Border color"

# TAB() - fully supported, just forgot to mark it before
$ws.Range("B165").Value = "X"
$ws.Range("C165").Value = "X"
$ws.Range("D165").Value = "X"
$ws.Range("E165").Value = "X"
